# Added headless flag, Test Summary and optimized Menu search flow
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "AUTO_OUTL24E5F"
$ws.Range("AY2").Value = "'03539463400"
$ws.Range("Q2").Value = 24.125126982276527

# Row 3
$ws.Range("B3").Value = "AUTO_OUTLF9BA1"
$ws.Range("AY3").Value = "'03540430200"
$ws.Range("Q3").Value = 24.349812618357568

# Row 4
$ws.Range("B4").Value = "AUTO_OUTL8A7D9"
$ws.Range("AY4").Value = "'03541077200"
$ws.Range("Q4").Value = 24.735438140229938
